# "real feed memes added"
#
# Two "meme" rows (old rows 66 & 67 of the "Numbers / Likes / Dislikes" feed
# table) get folded back into the top of the list: their Likes/Dislikes
# numbers are moved up into row 5 and row 22, and the now-empty rows 66/67
# are removed entirely. Finally the active selection is moved to F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 67 (Likes=34, Dislikes=6) -> moved into row 5
$ws.Range("B5").Value = 34
$ws.Range("C5").Value = 6

# Row 66 (Likes=29, Dislikes=11) -> moved into row 22
$ws.Range("B22").Value = 29
$ws.Range("C22").Value = 11

# The source rows (66 & 67) are now empty -- clear them out completely
# (values + formatting) so the rows disappear from the sheet, same as the
# other never-used trailing rows.
$ws.Range("A66:C67").Clear()

# Move the selection
$ws.Range("F5").Select()
